# Update of league bases (Estonia Meistriliiga), 14-05-2024 20:19
# 1) Some existing match rows had their data re-ordered/corrected (everything
#    except the running-index column A is swapped/rotated between rows).
# 2) Five new match rows (158-162) are appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AB (1-based column numbers 2..28) hold the per-match data;
# column A (1) is just the running match index and must stay untouched.
$firstCol = 2
$lastCol = 28

function Get-RowValues([int]$row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues([int]$row, $vals) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$i]
        $i++
    }
}

function Swap-Rows([int]$rowA, [int]$rowB) {
    $a = Get-RowValues $rowA
    $b = Get-RowValues $rowB
    Set-RowValues $rowA $b
    Set-RowValues $rowB $a
}

# --- Swap data of rows 64 and 65 ---
Swap-Rows 64 65

# --- Swap data of rows 95 and 96 ---
Swap-Rows 95 96

# --- Rotate data of rows 104 -> 106 -> 107 -> 104 (row 105 is unaffected) ---
$v104 = Get-RowValues 104
$v106 = Get-RowValues 106
$v107 = Get-RowValues 107
Set-RowValues 104 $v106
Set-RowValues 106 $v107
Set-RowValues 107 $v104

# --- Append 5 new match rows at the bottom (158..162) ---
# Use the last existing data row (157) as a formatting template so the new
# rows inherit the same cell styles (bold/bordered index column, date format).
$templateRow = 157

$newRows = @(
    @{ A=156; B=7721026; C="Estonia Meistriliiga"; D=45423.35416666666; E="FC Flora Tallinn";   F="Parnu JK Vaprus";      G=0; H=1; I="A"; J=1.285; K=5;   L=8;     M=1.285; N=5;   O=9;     P=-1.75; Q=1.95;  R=1.85;  S=3;    T=1.85;  U=1.95;  V=-1; W=-1;  X=8;     Y=-1;    Z=0.8500000000000001;  AA=-1; AB=0.95 },
    @{ A=157; B=7721027; C="Estonia Meistriliiga"; D=45423.45833333334; E="FC Kuressaare";      F="JK Tammeka Tartu";    G=0; H=0; I="D"; J=2.625; K=3.5; L=2.3;   M=2.4;   N=3.5; O=2.5;   P=0;     Q=1.85;  R=1.95;  S=2.75; T=1.85;  U=1.95;  V=-1; W=2.5;  X=-1;    Y=0;     Z=0;                    AA=-1; AB=0.95 },
    @{ A=158; B=7719655; C="Estonia Meistriliiga"; D=45424.35416666666; E="FC Levadia Tallinn"; F="JK Nomme Kalju";      G=0; H=0; I="D"; J=1.363; K=4.6; L=6.5;   M=1.5;   N=4.2; O=5.25;  P=-1.25; Q=1.95;  R=1.75;  S=2.75; T=1.85;  U=1.95;  V=-1; W=3.2;  X=-1;    Y=-1;    Z=0.75;                 AA=-1; AB=0.95 },
    @{ A=159; B=7719656; C="Estonia Meistriliiga"; D=45424.45833333334; E="JK Tallinna Kalev";  F="JK Trans Narva";      G=1; H=0; I="H"; J=1.75;  K=3.75;L=3.75;  M=2.2;   N=3.6; O=2.7;   P=-0.25; Q=2;     R=1.8;   S=2.75; T=1.825; U=1.975; V=1.2; W=-1;  X=-1;    Y=1;     Z=-1;                   AA=-1; AB=0.9750000000000001 },
    @{ A=160; B=7721028; C="Estonia Meistriliiga"; D=45424.54166666666; E="JK Nomme United";    F="Paide Linnameeskond"; G=0; H=1; I="A"; J=7;     K=4.8; L=1.333; M=6.5;   N=5;   O=1.333; P=1.5;   Q=1.925; R=1.875; S=3.25; T=2;     U=1.8;   V=-1; W=-1;  X=0.333; Y=0.925; Z=-1;                   AA=-1; AB=0.8 }
)

$rowIdx = 158
foreach ($row in $newRows) {
    # Copy formatting (styles) from the template row (only columns A:AB) first,
    # so the new row's cells (A: bold+border, D: datetime format, rest: default)
    # match the existing table's look without touching unrelated columns.
    $ws.Range("A" + $templateRow + ":AB" + $templateRow).Copy() | Out-Null
    $ws.Range("A" + $rowIdx + ":AB" + $rowIdx).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Cells.Item($rowIdx, 1).Value2  = $row.A
    $ws.Cells.Item($rowIdx, 2).Value2  = $row.B
    $ws.Cells.Item($rowIdx, 3).Value2  = $row.C
    $ws.Cells.Item($rowIdx, 4).Value2  = $row.D
    $ws.Cells.Item($rowIdx, 5).Value2  = $row.E
    $ws.Cells.Item($rowIdx, 6).Value2  = $row.F
    $ws.Cells.Item($rowIdx, 7).Value2  = $row.G
    $ws.Cells.Item($rowIdx, 8).Value2  = $row.H
    $ws.Cells.Item($rowIdx, 9).Value2  = $row.I
    $ws.Cells.Item($rowIdx, 10).Value2 = $row.J
    $ws.Cells.Item($rowIdx, 11).Value2 = $row.K
    $ws.Cells.Item($rowIdx, 12).Value2 = $row.L
    $ws.Cells.Item($rowIdx, 13).Value2 = $row.M
    $ws.Cells.Item($rowIdx, 14).Value2 = $row.N
    $ws.Cells.Item($rowIdx, 15).Value2 = $row.O
    $ws.Cells.Item($rowIdx, 16).Value2 = $row.P
    $ws.Cells.Item($rowIdx, 17).Value2 = $row.Q
    $ws.Cells.Item($rowIdx, 18).Value2 = $row.R
    $ws.Cells.Item($rowIdx, 19).Value2 = $row.S
    $ws.Cells.Item($rowIdx, 20).Value2 = $row.T
    $ws.Cells.Item($rowIdx, 21).Value2 = $row.U
    $ws.Cells.Item($rowIdx, 22).Value2 = $row.V
    $ws.Cells.Item($rowIdx, 23).Value2 = $row.W
    $ws.Cells.Item($rowIdx, 24).Value2 = $row.X
    $ws.Cells.Item($rowIdx, 25).Value2 = $row.Y
    $ws.Cells.Item($rowIdx, 26).Value2 = $row.Z
    $ws.Cells.Item($rowIdx, 27).Value2 = $row.AA
    $ws.Cells.Item($rowIdx, 28).Value2 = $row.AB

    $rowIdx++
}
